$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.524.66'
$ws.Range('E2').Value = '  -3.76%  '
$ws.Range('D3').Value = '3.575.38'
$ws.Range('E3').Value = '  -4.23%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '187.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('D7').Value = '3.572.28'
$ws.Range('E7').Value = '  -4.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.615'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.998'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.672'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.67%  '
$ws.Range('E11').Value = '  -9.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.91'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000262'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -10.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.80'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.12%  '
$ws.Range('D15').Value = '4.144.68'
$ws.Range('E15').Value = '  -4.28%  '
$ws.Range('D16').Value = '3.577.27'
$ws.Range('E16').Value = '  -4.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.58%  '
$ws.Range('D20').Value = '66.466.90'
$ws.Range('E20').Value = '  -3.62%  '
$ws.Range('E21').Value = '  -6.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '397.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.31%  '
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('E33').Value = '  -4.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '617.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.114'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.46%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '63.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.73%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.391'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.71%  '
$ws.Range('D40').Value = '0.0₃0760'
$ws.Range('E40').Value = '  -12.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.132'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Value = '3.026.49'
$ws.Range('E43').Value = '  +6.67%  '
$ws.Range('E44').Value = '  -7.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0409'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.96%  '
$ws.Range('E47').Value = '  -6.59%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.78%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '138.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.43%  '
$ws.Range('E51').Value = '  -0.12%  '
